# Actualizacion ejecuciones de casos de pruebas
# Updates the "Pasos" (steps) worksheet of the test-case workbook: rewords several
# steps/expected-results to use the <Ciudad1> placeholder, inserts a new step row
# ("Presionar el boton Buscar" / "Se abre un panel para configurar las estadisticas"),
# and renumbers/re-heights the remaining rows accordingly. Also updates the active
# sheet/selection to reflect that "Pasos" is now the sheet being worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pasos")

# --- Reword the first four existing steps -------------------------------------------------
$ws.Range("B2").Value = "Ingresar al sitio web."
$ws.Range("B3").Value = "Ingresar a la opcion Estadisticas"

$ws.Range("B4").Value = "Escribir <Ciudad1> en donde se debe ingresar la ciudad a buscar."
$ws.Range("C4").Value = "Se muestra el autocomplete con las opciones de las ciudad que se corresponden con <Ciudad1>"

$ws.Range("B5").Value = "Seleccionar <Ciudad1> en la lista de ciudad del autocomplete y presiono BUSCAR"
$ws.Range("C5").ClearContents()

# --- Insert the new step row between old row 5 and old row 6 ------------------------------
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = $null
$ws.Range("B6").Value = "Presionar el boton Buscar"
$ws.Range("C6").Value = "Se abre un panel para configurar las estadisticas"
$ws.Rows.Item(6).RowHeight = 15.75

# --- Reword the remaining (shifted-down) steps --------------------------------------------
$ws.Range("B7").Value = 'Presionar el boton "Nueva"'
$ws.Range("C7").Value = 'Se despliegan 2 opciones, "Historica" y "Tiempo Real"'

$ws.Range("B8").Value = 'Presionar el boton "Historica"'
$ws.Range("C8").Value = "Se muestra un panel para buscar estadisticas historicas"

$ws.Range("B9").Value = 'Ingresar Estaddistica: "Consulta", Por: "Playa",Desde 01/01/2014, Hasta 01/03/2015'

$ws.Range("B10").Value = 'Presionar el boton "Buscar"'
$ws.Range("C10").Value = "Se despliega el panel de estadisticas, con los filtros y las opciones del grafico"

# --- Row heights to match the reflowed (now mostly 2-line) step text ----------------------
$ws.Rows.Item(4).RowHeight = 26.25
$ws.Rows.Item(5).RowHeight = 26.25
$ws.Rows.Item(7).RowHeight = 26.25
$ws.Rows.Item(8).RowHeight = 26.25
$ws.Rows.Item(9).RowHeight = 26.25
$ws.Rows.Item(10).RowHeight = 26.25

# --- Selection / active sheet bookkeeping --------------------------------------------------
$ws.Range("B2:C10").Select()
$ws.Activate()

$datos = $wb.Worksheets.Item("DatosGenerales")
$datos.Range("A1").Select()
